# Reformat the single "questions = [...]" shared string from a Python
# dict/list repr (single-quoted) into pretty-printed JSON (indent=4,
# ASCII-escaped), and move it from A2 up into A1 - replacing the old
# placeholder value (0) that used to live in A1 with its bold/bordered
# style. The row that used to hold that placeholder is removed so the
# sheet ends up with just the single A1 cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newText = @'
questions = [
    {
        "title": "You are designing a new social media app. One of the app\u2019s features is the ability for users to post comments on other users' posts. You need to decompose the task of posting a comment into smaller steps.Which of the following is the best algorithm for suggesting search results?",
        "ques_type": 2,
        "options": [
            "1. Enter the text of the comment.2. Click the \"Post Comment\" button.3. Wait for the comment to be posted.",
            "1. Authenticate the user.2. Validate the comment text.3. Insert the comment into the database.",
            "1. Render the comment form.2. Handle the comment submission.3. Update the user's post with the new comment.",
            "1. Parse the comment text.2. Generate a preview of the comment.3. Post the comment to the user's feed."
        ],
        "score": "1. Authenticate the user.2. Validate the comment text.3. Insert the comment into the database."
    },
    {
        "title": "Which of the following would best represent the patterns in this list?x x y y y y z z x x x x",
        "ques_type": 2,
        "options": [
            "(x,2)(x,1)(z,1),(z,3)(x,1)(x,3)",
            "(x,2)(y,4)(z,2)(x,4)",
            "(x,3)(y,1)(z,3)",
            "(x,1)(y,4)(x,4)"
        ],
        "score": "(x,2)(y,4)(z,2)(x,4)"
    },
    {
        "title": "You are a delivery driver for a package delivery company. You have a list of packages to deliver to different addresses in the city. You want to spend as little time as possible driving around the city delivering the packages. Which of the following is the best algorithm for doing this?",
        "ques_type": 2,
        "options": [
            "Drive to each address on the list in order. If there is a package to deliver at that address, then deliver it.",
            "For each package on the list, go to the address it is to be delivered to and deliver it.",
            "Determine which neighborhood each delivery address is in. Cluster the list by deliveries in the same neighborhood. Drive to each neighborhood once and deliver all of the packages in that neighborhood.",
            "Drive to the address of the package that is furthest away from your starting location. Deliver that package. Then drive to the address of the next furthest away package, and so on until all of the packages have been delivered."
        ],
        "score": "Determine which neighborhood each delivery address is in. Cluster the list by deliveries in the same neighborhood. Drive to each neighborhood once and deliver all of the packages in that neighborhood."
    },
    {
        "title": "You are designing a software application for a library management system. The system needs to keep track of various details for each book, such as title, author, publication date, and ISBN. However, you also need to ensure the software can handle other types of items, like DVDs and magazines, with their own specific attributes.Which concept of abstraction best applies to this situation?",
        "ques_type": 2,
        "options": [
            "Reducing the storage space required for storing book details in the database.",
            "Simplifying the user interface to make it more user-friendly.",
            "Generalizing the common attributes and behaviors of different types of library items.",
            "Increasing the number of available functions and features in the software."
        ],
        "score": "Generalizing the common attributes and behaviors of different types of library items."
    }
]
'@

# A2 (the cell that used to hold the long string) is no longer needed -
# clear both its content and formatting so the used range shrinks back
# down to a single cell.
$ws.Range("A2").Clear()

# A1 used to hold the number 0 with a bold/centered/bordered style;
# strip that formatting before writing the new text so A1 ends up with
# plain/default formatting.
$ws.Range("A1").ClearFormats()
$ws.Range("A1").Value = $newText

# Writing a multi-line value auto-expands the row height; put it back to
# the sheet's default (no explicit/custom row height).
$ws.Rows(1).EntireRow.AutoFit()
